$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.133.38"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.815.52"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'706.44"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "'171.65"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").Value = "3.814.84"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'7.71"
$ws.Range("E11").Value = "  +7.06%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "'35.85"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "4.454.44"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "3.832.12"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "71.086.26"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "'17.49"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'501.62"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'84.25"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").Value = "3.962.74"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "'10.36"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").Value = "'2.27"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "'7.36"
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").Value = "'29.08"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").Value = "'0.174"
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").Value = "3.777.63"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'9.10"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("B38").Value = "Binance-PegBSC-USD"
$ws.Range("C38").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").Value = "'2.37"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "'167.14"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "'0.000315"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "'49.10"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'421.01"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'0.294"
$ws.Range("E51").Value = "  -2.21%  "
